$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching the formatting of the
# neighboring header cell (G1) by copying its format over.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the corresponding value in H2 (plain, unstyled numeric cell like B2:G2)
$ws.Range("H2").Value = 0
